$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells stay formatted as Text so that
# numeric-looking values (e.g. 569.16) are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.659.04'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = '3.409.74'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '569.16'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').Value = '157.32'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.413.12'
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('D9').Value = '0.567'
$ws.Range('E9').Value = '  -8.52%  '
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('E11').Value = '  -4.27%  '
$ws.Range('D12').Value = '0.422'
$ws.Range('E12').Value = '  -4.65%  '
$ws.Range('D13').Value = '3.994.47'
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('D15').Value = '26.98'
$ws.Range('E15').Value = '  -3.41%  '
$ws.Range('D16').Value = '0.0000172'
$ws.Range('E16').Value = '  -8.73%  '
$ws.Range('D17').Value = '63.718.56'
$ws.Range('E17').Value = '  -1.39%  '
$ws.Range('D18').Value = '3.414.30'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('E19').Value = '  -4.65%  '
$ws.Range('D20').Value = '13.57'
$ws.Range('D21').Value = '385.90'
$ws.Range('E21').Value = '  +1.88%  '
$ws.Range('D22').Value = '7.75'
$ws.Range('E22').Value = '  -3.59%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').Value = '71.18'
$ws.Range('E24').Value = '  -1.84%  '
$ws.Range('D25').Value = '0.515'
$ws.Range('E25').Value = '  -6.65%  '
$ws.Range('E26').Value = '  -4.52%  '
$ws.Range('E27').Value = '  -6.28%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = '6.08'
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('D31').Value = '1.39'
$ws.Range('E31').Value = '  -6.84%  '
$ws.Range('E32').Value = '  -2.74%  '
$ws.Range('D34').Value = '22.85'
$ws.Range('E34').Value = '  -1.34%  '
$ws.Range('E35').Value = '  -4.36%  '
$ws.Range('D36').Value = '1.51'
$ws.Range('E36').Value = '  -6.70%  '
$ws.Range('D37').Value = '160.64'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').Value = '0.842'
$ws.Range('E38').Value = '  +8.79%  '
$ws.Range('D39').Value = '1.82'
$ws.Range('E39').Value = '  -4.73%  '
$ws.Range('D40').Value = '2.785.32'
$ws.Range('E40').Value = '  -3.25%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '0.0723'
$ws.Range('E41').Value = '  -5.49%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '25.81'
$ws.Range('E42').Value = '  -4.54%  '
$ws.Range('D43').Value = '42.99'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '25.90'
$ws.Range('E44').Value = '  -2.72%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '6.35'
$ws.Range('E45').Value = '  -8.95%  '
$ws.Range('D46').Value = '4.34'
$ws.Range('E46').Value = '  -5.97%  '
$ws.Range('E47').Value = '  -4.91%  '
$ws.Range('D48').Value = '2.37'
$ws.Range('E48').Value = '  +8.49%  '
$ws.Range('D49').Value = '326.64'
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('E50').Value = '  -5.39%  '
$ws.Range('E51').Value = '  -5.03%  '
